$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 44 (Excel shifts rows 44-60
# down to 45-61, preserving all their existing values/styles).
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new market record.
$ws.Range("A44").Value = 11
$ws.Range("B44").Value = "Vega Monumental Concepción"
$ws.Range("C44").Value = "Bíobío"
$ws.Range("D44").Value = 45007
$ws.Range("E44").Value = 8
$ws.Range("F44").Value = 100112043
$ws.Range("G44").Value = "Pepino dulce"
$ws.Range("H44").Value = "Cultivar IV Región"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 100
$ws.Range("K44").Value = 13000
$ws.Range("L44").Value = 14000
$ws.Range("M44").Value = 13500
$ws.Range("N44").Value = "$/bandeja 18 kilos"
$ws.Range("O44").Value = "Provincia de Limarí"
$ws.Range("P44").Value = 750
$ws.Range("Q44").Value = 18
$ws.Range("R44").Value = "Hortaliza"

# Note: Rows.Item(44).Insert() already carries the "Fecha" column's
# date style (numFmt "YYYY-MM-DD HH:MM:SS", style index 2) down into the
# newly inserted row, matching the other rows, so no extra style copy is
# needed here.
